$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set explicit column widths for A:D (as previously auto-fitted by the author)
$ws.Columns("A").ColumnWidth = 13.666666666666666
$ws.Columns("B").ColumnWidth = 16.166666666666668
$ws.Columns("C").ColumnWidth = 18.166666666666668
$ws.Columns("D").ColumnWidth = 17.666666666666668

# Update the current selection/active cell shown when the sheet is reopened
$ws.Range("E37").Select() | Out-Null

# Best-effort: reflect the new window size/position captured in the saved view
$win = $excel.ActiveWindow
$win.Top = 465
$win.Left = 3690
$win.Width = 21600
$win.Height = 13560
